# Add two new columns, I (header "I0") and J (header "IF"), to Sheet1.
# Column H already has header "IP" and style index 1 (bold/bordered/centered).
# We copy that header style onto the new header cells, then fill in the
# per-row numeric data for rows 2-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells (row 1), matching the style used by the existing header cells.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats: copy H1's style (bold/border/centered)

# Data values for rows 2-43: column I ("I0") and column J ("IF").
$dataI = @(8,5,6,1,5,7,4,7,6,6,6,8,1,4,1,1,6,9,7,6,8,9,7,7,7,7,7,1,1,6,6,10,8,7,9,6,7,4,6,5,6,6)
$dataJ = @(9,5,6,2,5,7,5,7,6,7,6,9,1,4,2,3,6,9,8,7,8,9,9,7,8,8,8,1,1,6,6,12,8,8,9,7,7,5,7,6,6,6)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
